$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 0.5797639999999999
$ws.Range("N2").Value = 0.691148
$ws.Range("O2").Value = 0.04658545143903391
$ws.Range("P2").Value = 0.04658545143903391
$ws.Range("Q2").Value = 0.04452252545244444
$ws.Range("R2").Value = 0.4007027290719999
$ws.Range("S2").Value = 0.04658545143903391
$ws.Range("T2").Value = 0.04658545143903391

# Row 3
$ws.Range("H3").Value = 0.5797639999999999
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.8771473333333333
$ws.Range("N3").Value = 2.631442
$ws.Range("O3").Value = 0.177367095767671
$ws.Range("P3").Value = 0.177367095767671
$ws.Range("Q3").Value = 0.1695128155208889
$ws.Range("R3").Value = 1.525615339688
$ws.Range("S3").Value = 0.177367095767671
$ws.Range("T3").Value = 0.177367095767671

# Row 4
$ws.Range("H4").Value = 0.5797639999999999
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.09376766666666668
$ws.Range("N4").Value = 0.281303
$ws.Range("O4").Value = 0.01896066724660212
$ws.Range("P4").Value = 0.01896066724660211
$ws.Range("Q4").Value = 0.01812103916577778
$ws.Range("R4").Value = 0.163089352492
$ws.Range("S4").Value = 0.01896066724660212
$ws.Range("T4").Value = 0.01896066724660211

# Row 5
$ws.Range("H5").Value = 0.5797639999999999
$ws.Range("M5").Value = 2.620343
$ws.Range("N5").Value = 7.861029
$ws.Range("O5").Value = 0.5298569694773585
$ws.Range("P5").Value = 0.5298569694773585
$ws.Range("Q5").Value = 0.5063935130173334
$ws.Range("R5").Value = 4.557541617156
$ws.Range("S5").Value = 0.5298569694773585
$ws.Range("T5").Value = 0.5298569694773585

# Row 6
$ws.Range("H6").Value = 0.5797639999999999
$ws.Range("M6").Value = 1.123737333333333
$ws.Range("N6").Value = 3.371212
$ws.Range("O6").Value = 0.2272298160693345
$ws.Range("P6").Value = 0.2272298160693345
$ws.Range("Q6").Value = 0.2171674837742222
$ws.Range("R6").Value = 1.954507353968
$ws.Range("S6").Value = 0.2272298160693345
$ws.Range("T6").Value = 0.2272298160693345
